# Update PriceHistory data:
# Insert 12 new (most recent) price rows at the top of the data table,
# just below the header row, pushing the existing history down by 12 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PriceHistory")

# Insert 12 blank rows right after the header (row 1)
$ws.Range("A2:A13").EntireRow.Insert()

# Make sure the new cells are stored as text (matching the rest of the
# Price Date / Price / Currency columns, which are all text-formatted).
$newRange = $ws.Range("A2:C13")
$newRange.NumberFormat = "@"

$dates    = @("04/11/2025","03/11/2025","31/10/2025","30/10/2025","29/10/2025","28/10/2025","27/10/2025","24/10/2025","23/10/2025","22/10/2025","21/10/2025","17/10/2025")
$prices   = @("1.057","1.063","1.055","1.058","1.062","1.051","1.054","1.054","1.045","1.046","1.052","1.024")
$currency = "SGD"

for ($i = 0; $i -lt $dates.Length; $i++) {
    $r = 2 + $i
    $ws.Cells.Item($r, 1).Value = $dates[$i]
    $ws.Cells.Item($r, 2).Value = $prices[$i]
    $ws.Cells.Item($r, 3).Value = $currency
}
